# Insert a new data row just before the current row 339 (Macroferia Regional
# de Talca / Coliflor, weekly price series). This shifts the existing rows
# 339-442 down to 340-443 (preserving all of their data), and grows the
# sheet's used range from A1:R442 to A1:R443. The freshly inserted row 339
# is then populated with the new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 339:442 down to 340:443, leaving row 339 empty (but formatted
# like its neighbours, since Excel's Insert copies formatting from above).
$ws.Rows("339").Insert()

# Fill in the newly-inserted row with the new observation.
$ws.Range("A339").Value = 5
$ws.Range("B339").Value = "Macroferia Regional de Talca"
$ws.Range("C339").Value = "Maule"
$ws.Range("D339").Value = 45093
$ws.Range("E339").Value = 7
$ws.Range("F339").Value = 100112008
$ws.Range("G339").Value = "Coliflor"
$ws.Range("H339").Value = "Sin especificar"
$ws.Range("I339").Value = "Primera"
$ws.Range("J339").Value = 6000
$ws.Range("K339").Value = 600
$ws.Range("L339").Value = 700
$ws.Range("M339").Value = 650
$ws.Range("N339").Value = "$/unidad"
$ws.Range("O339").Value = "Región del Maule"
$ws.Range("P339").Value = 650
$ws.Range("Q339").Value = 1
$ws.Range("R339").Value = "Hortaliza"
